$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data update: row 2 "usuario" test value changed from "testing10" to "pruebauser01"
$ws.Range("D2").Value = "pruebauser01"

# Selection/view state: active cell moved to D7 (also resets the scrolled
# top-left cell back into view)
$ws.Range("D7").Select()
